$wb = $excel.ActiveWorkbook

# --- Sheet "data" ---
$ws1 = $wb.Worksheets.Item("data")

# Header date fix: 25. 1. 2021 -> 25. 1. 2022
$ws1.Range("H1").Value = "25. 1. 2022"

# Row 380 - "Jiz ockovan/a treti davkou"
$ws1.Range("E380").Value = 0.04
$ws1.Range("F380").Value = 0.07000000000000001
$ws1.Range("G380").Value = 0.14

# Row 381 - "Rozhodne ano"
$ws1.Range("D381").Value = 0.25
$ws1.Range("E381").Value = 0.23
$ws1.Range("F381").Value = 0.32
$ws1.Range("G381").Value = 0.22

# Row 382 - "Spise ano"
$ws1.Range("D382").Value = 0.19
$ws1.Range("E382").Value = 0.15
$ws1.Range("F382").Value = 0.17

# Row 383 - "Nevim"
$ws1.Range("D383").Value = 0.09
$ws1.Range("F383").Value = 0.04

# Row 384 - "Spise ne"
$ws1.Range("E384").Value = 0.09
$ws1.Range("F384").Value = 0.07000000000000001

# Row 385 - "Rozhodne ne"
$ws1.Range("D385").Value = 0.08
$ws1.Range("E385").Value = 0.05

# Row 386 - "Neni ockovan/a ani jednou davkou"
$ws1.Range("D386").Value = 0.32
$ws1.Range("E386").Value = 0.32

# Row 387 - new data for "Bez nakazy koronavirem" / "Jiz ockovan/a treti davkou"
$ws1.Range("D387").Value = 0
$ws1.Range("E387").Value = 0.04
$ws1.Range("F387").Value = 0.11
$ws1.Range("G387").Value = 0.3
$ws1.Range("H387").Value = 0.44

# Row 388 - "Rozhodne ano"
$ws1.Range("D388").Value = 0.25
$ws1.Range("E388").Value = 0.33
$ws1.Range("F388").Value = 0.37
$ws1.Range("G388").Value = 0.21
$ws1.Range("H388").Value = 0.09

# Row 389 - "Spise ano"
$ws1.Range("D389").Value = 0.21
$ws1.Range("E389").Value = 0.17
$ws1.Range("F389").Value = 0.17
$ws1.Range("G389").Value = 0.13
$ws1.Range("H389").Value = 0.11

# Row 390 - "Nevim"
$ws1.Range("D390").Value = 0.12
$ws1.Range("E390").Value = 0.09
$ws1.Range("F390").Value = 0.07000000000000001
$ws1.Range("G390").Value = 0.07000000000000001
$ws1.Range("H390").Value = 0.06

# Row 391 - "Spise ne"
$ws1.Range("D391").Value = 0.06
$ws1.Range("E391").Value = 0.06
$ws1.Range("F391").Value = 0.05
$ws1.Range("G391").Value = 0.04
$ws1.Range("H391").Value = 0.06

# Row 392 - "Rozhodne ne"
$ws1.Range("D392").Value = 0.04
$ws1.Range("E392").Value = 0.03
$ws1.Range("F392").Value = 0.02
$ws1.Range("G392").Value = 0.02
$ws1.Range("H392").Value = 0.03

# Row 393 - "Neni ockovan/a ani jednou davkou"
$ws1.Range("D393").Value = 0.32
$ws1.Range("E393").Value = 0.28
$ws1.Range("F393").Value = 0.21
$ws1.Range("G393").Value = 0.23
$ws1.Range("H393").Value = 0.21

# --- Sheet "pocetR" ---
$ws2 = $wb.Worksheets.Item("pocetR")

# Header date fix: 25. 1. 2021 -> 25. 1. 2022
$ws2.Range("G1").Value = "25. 1. 2022"

# Row 57 - "V minulosti byla testem potvrzena nakaza koronavirem"
$ws2.Range("C57").Value = 262
$ws2.Range("D57").Value = 273
$ws2.Range("E57").Value = 322
$ws2.Range("F57").Value = 366

# Row 58 - "Bez nakazy koronavirem (potvrzene testem)"
$ws2.Range("C58").Value = 1593
$ws2.Range("D58").Value = 1436
$ws2.Range("E58").Value = 1468
$ws2.Range("F58").Value = 1401
